$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.481
$ws.Range("F4").Value = 0.053
$ws.Range("G4").Value = 0.23
$ws.Range("N4").Value = 0.485
$ws.Range("O4").Value = 0.057
$ws.Range("P4").Value = 0.24
$ws.Range("Q4").Value = 0.052
$ws.Range("R4").Value = 0.035
$ws.Range("S4").Value = 0.188
$ws.Range("W4").Value = 0.367
$ws.Range("AI4").Value = 0.403
$ws.Range("AJ4").Value = 0.096
$ws.Range("AK4").Value = 0.309
$ws.Range("AU4").Value = 0.24
$ws.Range("AW4").Value = 0.163
$ws.Range("BA4").Value = 2.041
$ws.Range("BB4").Value = 0.142
$ws.Range("BC4").Value = 0.377
$ws.Range("BG4").Value = 0.722
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.75
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.68
$ws.Range("BQ4").Value = 0.761
$ws.Range("E5").Value = 0.606
$ws.Range("F5").Value = 0.058
$ws.Range("G5").Value = 0.241
$ws.Range("N5").Value = 0.733
$ws.Range("O5").Value = 0.065
$ws.Range("P5").Value = 0.255
$ws.Range("Q5").Value = 0.035
$ws.Range("R5").Value = 0.014
$ws.Range("S5").Value = 0.12
$ws.Range("W5").Value = 0.335
$ws.Range("X5").Value = 0.1
$ws.Range("Y5").Value = 0.316
$ws.Range("AI5").Value = 0.404
$ws.Range("AJ5").Value = 0.092
$ws.Range("AK5").Value = 0.303
$ws.Range("AU5").Value = 0.448
$ws.Range("AV5").Value = 0.079
$ws.Range("AW5").Value = 0.28
$ws.Range("BA5").Value = 1.306
$ws.Range("BB5").Value = 0.074
$ws.Range("BC5").Value = 0.271
$ws.Range("BG5").Value = 0.383
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.225
$ws.Range("BM5").Value = 0.525
$ws.Range("BN5").Value = 0.047
$ws.Range("BO5").Value = 0.217
$ws.Range("BP5").Value = 0.435
$ws.Range("BQ5").Value = 0.458
$ws.Range("E6").Value = 0.536
$ws.Range("N6").Value = 0.584
$ws.Range("Q6").Value = 0.042
$ws.Range("W6").Value = 0.35
$ws.Range("AI6").Value = 0.403
$ws.Range("AU6").Value = 0.313
$ws.Range("BA6").Value = 1.586
$ws.Range("BG6").Value = 0.5
$ws.Range("BM6").Value = 0.618
$ws.Range("BP6").Value = 0.529
$ws.Range("BQ6").Value = 0.569
$ws.Range("E7").Value = 0.576
$ws.Range("N7").Value = 0.665
$ws.Range("Q7").Value = 0.037
$ws.Range("W7").Value = 0.341
$ws.Range("AI7").Value = 0.404
$ws.Range("AU7").Value = 0.382
$ws.Range("BA7").Value = 1.405
$ws.Range("BG7").Value = 0.423
$ws.Range("BM7").Value = 0.5590000000000001
$ws.Range("BP7").Value = 0.468
$ws.Range("BQ7").Value = 0.497
$ws.Range("E8").Value = 0.705
$ws.Range("F8").Value = 0.07199999999999999
$ws.Range("G8").Value = 0.268
$ws.Range("N8").Value = 0.821
$ws.Range("O8").Value = 0.044
$ws.Range("P8").Value = 0.21
$ws.Range("Q8").Value = 0.038
$ws.Range("W8").Value = 0.407
$ws.Range("X8").Value = 0.122
$ws.Range("Y8").Value = 0.349
$ws.Range("AI8").Value = 0.472
$ws.Range("AJ8").Value = 0.14
$ws.Range("AK8").Value = 0.374
$ws.Range("AU8").Value = 0.392
$ws.Range("AV8").Value = 0.08500000000000001
$ws.Range("AW8").Value = 0.292
$ws.Range("BA8").Value = 1.77
$ws.Range("BB8").Value = 0.109
$ws.Range("BC8").Value = 0.33
$ws.Range("BG8").Value = 0.57
$ws.Range("BH8").Value = 0.11
$ws.Range("BI8").Value = 0.331
$ws.Range("BM8").Value = 0.673
$ws.Range("BN8").Value = 0.062
$ws.Range("BO8").Value = 0.249
$ws.Range("BP8").Value = 0.59
$ws.Range("BQ8").Value = 0.626
$ws.Range("E9").Value = 0.667
$ws.Range("F9").Value = 0.222
$ws.Range("G9").Value = 0.471
$ws.Range("N9").Value = 0.762
$ws.Range("O9").Value = 0.181
$ws.Range("P9").Value = 0.426
$ws.Range("W9").Value = 0.31
$ws.Range("X9").Value = 0.214
$ws.Range("Y9").Value = 0.462
$ws.Range("AI9").Value = 0.429
$ws.Range("AJ9").Value = 0.245
$ws.Range("AK9").Value = 0.495
$ws.Range("BA9").Value = 1.738
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.619
$ws.Range("BH9").Value = 0.236
$ws.Range("BI9").Value = 0.486
$ws.Range("BM9").Value = 0.643
$ws.Range("BN9").Value = 0.23
$ws.Range("BO9").Value = 0.479
$ws.Range("BP9").Value = 0.579
$ws.Range("BQ9").Value = 0.618
$ws.Range("E10").Value = 0.8100000000000001
$ws.Range("F10").Value = 0.154
$ws.Range("G10").Value = 0.393
$ws.Range("N10").Value = 0.952
$ws.Range("O10").Value = 0.045
$ws.Range("P10").Value = 0.213
$ws.Range("W10").Value = 0.524
$ws.Range("X10").Value = 0.249
$ws.Range("Y10").Value = 0.499
$ws.Range("AI10").Value = 0.5
$ws.Range("AJ10").Value = 0.25
$ws.Range("AK10").Value = 0.5
$ws.Range("AU10").Value = 0.381
$ws.Range("AV10").Value = 0.236
$ws.Range("AW10").Value = 0.486
$ws.Range("BA10").Value = 2.19
$ws.Range("BB10").Value = 0.214
$ws.Range("BC10").Value = 0.462
$ws.Range("BG10").Value = 0.6899999999999999
$ws.Range("BH10").Value = 0.214
$ws.Range("BI10").Value = 0.462
$ws.Range("BM10").Value = 0.8100000000000001
$ws.Range("BN10").Value = 0.154
$ws.Range("BO10").Value = 0.393
$ws.Range("BP10").Value = 0.73
$ws.Range("BQ10").Value = 0.764
$ws.Range("E11").Value = 0.857
$ws.Range("F11").Value = 0.122
$ws.Range("G11").Value = 0.35
$ws.Range("N11").Value = 0.952
$ws.Range("O11").Value = 0.045
$ws.Range("P11").Value = 0.213
$ws.Range("W11").Value = 0.524
$ws.Range("X11").Value = 0.249
$ws.Range("Y11").Value = 0.499
$ws.Range("AI11").Value = 0.571
$ws.Range("AJ11").Value = 0.245
$ws.Range("AK11").Value = 0.495
$ws.Range("AU11").Value = 0.548
$ws.Range("AV11").Value = 0.248
$ws.Range("AW11").Value = 0.498
$ws.Range("BA11").Value = 2.19
$ws.Range("BB11").Value = 0.214
$ws.Range("BC11").Value = 0.462
$ws.Range("BG11").Value = 0.6899999999999999
$ws.Range("BH11").Value = 0.214
$ws.Range("BI11").Value = 0.462
$ws.Range("BM11").Value = 0.8100000000000001
$ws.Range("BN11").Value = 0.154
$ws.Range("BO11").Value = 0.393
$ws.Range("BP11").Value = 0.73
$ws.Range("BQ11").Value = 0.77
$ws.Range("E12").Value = 1.417
$ws.Range("F12").Value = 0.854
$ws.Range("G12").Value = 0.924
$ws.Range("N12").Value = 1.25
$ws.Range("O12").Value = 0.287
$ws.Range("P12").Value = 0.536
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.432
$ws.Range("Y12").Value = 0.657
$ws.Range("AI12").Value = 1.583
$ws.Range("AJ12").Value = 1.493
$ws.Range("AK12").Value = 1.222
$ws.Range("AU12").Value = 2.88
$ws.Range("AV12").Value = 3.466
$ws.Range("AW12").Value = 1.862
$ws.Range("BA12").Value = 3.787
$ws.Range("BB12").Value = 0.449
$ws.Range("BC12").Value = 0.67
$ws.Range("BG12").Value = 1.138
$ws.Range("BH12").Value = 0.188
$ws.Range("BI12").Value = 0.433
$ws.Range("BM12").Value = 1.235
$ws.Range("BN12").Value = 0.239
$ws.Range("BO12").Value = 0.489
$ws.Range("BP12").Value = 1.262
$ws.Range("BQ12").Value = 1.248
$ws.Range("E13").Value = 1.415
$ws.Range("F13").Value = 0.295
$ws.Range("G13").Value = 0.543
$ws.Range("N13").Value = 1.737
$ws.Range("O13").Value = 0.466
$ws.Range("P13").Value = 0.6830000000000001
$ws.Range("W13").Value = 0.985
$ws.Range("X13").Value = 0.199
$ws.Range("Y13").Value = 0.446
$ws.Range("AI13").Value = 1.159
$ws.Range("AJ13").Value = 0.312
$ws.Range("AK13").Value = 0.5580000000000001
$ws.Range("AU13").Value = 2.048
$ws.Range("AV13").Value = 0.344
$ws.Range("AW13").Value = 0.587
$ws.Range("BA13").Value = 2.187
$ws.Range("BB13").Value = 0.278
$ws.Range("BC13").Value = 0.527
$ws.Range("BG13").Value = 0.547
$ws.Range("BH13").Value = 0.05
$ws.Range("BI13").Value = 0.224
$ws.Range("BM13").Value = 0.787
$ws.Range("BN13").Value = 0.163
$ws.Range("BO13").Value = 0.403
$ws.Range("BP13").Value = 0.729
$ws.Range("BQ13").Value = 0.667
